$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps being treated as text, matching the
# workbook's existing data (values like "1.001" or "236.19" would
# otherwise be auto-converted to numbers by Excel's smart input parsing).
$ws.Range("D2:D51").NumberFormat = "@"

    $ws.Range("D2").Value = "30.546.25"
    $ws.Range("E2").Value = "  -0.37%  "
    $ws.Range("D3").Value = "1.876.65"
    $ws.Range("E3").Value = "  -0.75%  "
    $ws.Range("E4").Value = "  +0.02%  "
    $ws.Range("D5").Value = "236.19"
    $ws.Range("E5").Value = "  -3.53%  "
    $ws.Range("D6").Value = "1.001"
    $ws.Range("E6").Value = "  -0.02%  "
    $ws.Range("D7").Value = "0.4876"
    $ws.Range("E7").Value = "  -1.82%  "
    $ws.Range("D8").Value = "0.2896"
    $ws.Range("E8").Value = "  -2.11%  "
    $ws.Range("D9").Value = "0.06667"
    $ws.Range("E9").Value = "  -2.11%  "
    $ws.Range("D10").Value = "1.875.42"
    $ws.Range("E10").Value = "  -0.78%  "
    $ws.Range("D11").Value = "16.57"
    $ws.Range("E11").Value = "  -2.79%  "
    $ws.Range("D12").Value = "0.07239"
    $ws.Range("E12").Value = "  -0.95%  "
    $ws.Range("D13").Value = "88.61"
    $ws.Range("E13").Value = "  -2.45%  "
    $ws.Range("D14").Value = "4.995"
    $ws.Range("E14").Value = "  -1.43%  "
    $ws.Range("D15").Value = "0.6499"
    $ws.Range("E15").Value = "  -3.26%  "
    $ws.Range("D16").Value = "30.489.87"
    $ws.Range("E16").Value = "  -0.49%  "
    $ws.Range("D17").Value = "0.000007851"
    $ws.Range("E17").Value = "  -0.95%  "
    $ws.Range("D18").Value = "1.001"
    $ws.Range("E18").Value = "  -0.15%  "
    $ws.Range("D19").Value = "12.98"
    $ws.Range("E19").Value = "  -1.56%  "
    $ws.Range("D20").Value = "2.116.08"
    $ws.Range("E20").Value = "  -0.88%  "
    $ws.Range("D21").Value = "1.001"
    $ws.Range("E21").Value = "  +0.12%  "
    $ws.Range("D22").Value = "4.713"
    $ws.Range("E22").Value = "  -2.81%  "
    $ws.Range("D23").Value = "197.67"
    $ws.Range("E23").Value = "  +12.23%  "
    $ws.Range("D24").Value = "6.123"
    $ws.Range("E24").Value = "  +1.11%  "
    $ws.Range("D25").Value = "9.360"
    $ws.Range("E25").Value = "  +0.93%  "
    $ws.Range("D26").Value = "157.08"
    $ws.Range("E26").Value = "  +1.40%  "
    $ws.Range("D27").Value = "18.46"
    $ws.Range("E27").Value = "  -1.49%  "
    $ws.Range("D28").Value = "1.827"
    $ws.Range("E28").Value = "  -5.07%  "
    $ws.Range("D29").Value = "1.409"
    $ws.Range("E29").Value = "  +1.36%  "
    $ws.Range("D30").Value = "4.249"
    $ws.Range("E30").Value = "  -1.86%  "
    $ws.Range("D31").Value = "0.09019"
    $ws.Range("E31").Value = "  +1.15%  "
    $ws.Range("D32").Value = "3.918"
    $ws.Range("E32").Value = "  -2.71%  "
    $ws.Range("D33").Value = "0.05113"
    $ws.Range("E33").Value = "  -2.12%  "
    $ws.Range("D34").Value = "0.7209"
    $ws.Range("E34").Value = "  -2.45%  "
    $ws.Range("D35").Value = "1.078"
    $ws.Range("E35").Value = "  -5.00%  "
    $ws.Range("D36").Value = "2.691"
    $ws.Range("E36").Value = "  +0.30%  "
    $ws.Range("E37").Value = "  -3.11%  "
    $ws.Range("D38").Value = "2.663"
    $ws.Range("E38").Value = "  -1.49%  "
    $ws.Range("D39").Value = "0.9200"
    $ws.Range("E39").Value = "  -1.60%  "
    $ws.Range("D40").Value = "2.042"
    $ws.Range("E40").Value = "  -5.91%  "
    $ws.Range("D41").Value = "0.4384"
    $ws.Range("E41").Value = "  +0.54%  "
    $ws.Range("D42").Value = "104.93"
    $ws.Range("E42").Value = "  -0.56%  "
    $ws.Range("D43").Value = "0.9950"
    $ws.Range("E43").Value = "  -0.67%  "
    $ws.Range("D44").Value = "5.711"
    $ws.Range("E44").Value = "  -1.74%  "
    $ws.Range("D45").Value = "0.1327"
    $ws.Range("E45").Value = "  -2.15%  "
    $ws.Range("D46").Value = "7.361"
    $ws.Range("E46").Value = "  -3.85%  "
    $ws.Range("E47").Value = "  +3.72%  "
    $ws.Range("D48").Value = "0.05825"
    $ws.Range("E48").Value = "  -0.01%  "
    $ws.Range("D49").Value = "8.612"
    $ws.Range("E49").Value = "  +1.06%  "
    $ws.Range("D50").Value = "1.401"
    $ws.Range("E50").Value = "  +1.62%  "
    $ws.Range("D51").Value = "33.11"
    $ws.Range("E51").Value = "  -0.92%  "
